$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 3).Value = 0.3897807598114014
$ws.Cells.Item(2, 5).Value = 241.0084293445489
$ws.Cells.Item(2, 6).Value = 0.00814999997897378
$ws.Cells.Item(2, 7).Value = 0.006757097845831176
$ws.Cells.Item(2, 8).Value = 0.006255734650482539
$ws.Cells.Item(2, 9).Value = 0.005978640047716665
$ws.Cells.Item(2, 10).Value = 0.005684635904777876
$ws.Cells.Item(2, 11).Value = 0.005684575056065472
$ws.Cells.Item(2, 12).Value = 0.005598107917332675
$ws.Cells.Item(2, 13).Value = 0.005299269240548413
$ws.Cells.Item(2, 14).Value = 0.005267900616459141
$ws.Cells.Item(2, 15).Value = 0.005267900616459141
$ws.Cells.Item(2, 16).Value = 0.00525324927168444
$ws.Cells.Item(2, 17).Value = 0.005150371800978338
$ws.Cells.Item(2, 18).Value = 0.005032099166538388
$ws.Cells.Item(2, 19).Value = 0.004977470807231996
$ws.Cells.Item(2, 20).Value = 0.004858849401079368
$ws.Cells.Item(2, 21).Value = 0.00483122908839799
$ws.Cells.Item(2, 22).Value = 0.004772651970255197
$ws.Cells.Item(2, 23).Value = 0.004746535171009339
$ws.Cells.Item(2, 24).Value = 0.004726611766995576
$ws.Cells.Item(2, 25).Value = 0.004698020065195883
$ws.Cells.Item(3, 3).Value = 0.3749868869781494
$ws.Cells.Item(3, 5).Value = 243.1964959708384
$ws.Cells.Item(3, 6).Value = 0.008537041040833346
$ws.Cells.Item(3, 7).Value = 0.006541561275571047
$ws.Cells.Item(3, 8).Value = 0.006396265637847665
$ws.Cells.Item(3, 9).Value = 0.005930994332067452
$ws.Cells.Item(3, 10).Value = 0.005896673940110808
$ws.Cells.Item(3, 11).Value = 0.005518176093702497
$ws.Cells.Item(3, 12).Value = 0.005356391950690717
$ws.Cells.Item(3, 13).Value = 0.005356391950690717
$ws.Cells.Item(3, 14).Value = 0.005275702546874299
$ws.Cells.Item(3, 15).Value = 0.005275702546874299
$ws.Cells.Item(3, 16).Value = 0.005275702546874299
$ws.Cells.Item(3, 17).Value = 0.005109974219006521
$ws.Cells.Item(3, 18).Value = 0.005072359058979637
$ws.Cells.Item(3, 19).Value = 0.004960225145065587
$ws.Cells.Item(3, 20).Value = 0.00493687020071461
$ws.Cells.Item(3, 21).Value = 0.004896641494617308
$ws.Cells.Item(3, 22).Value = 0.004825035400355833
$ws.Cells.Item(3, 23).Value = 0.004801243518996363
$ws.Cells.Item(3, 24).Value = 0.004769929660607041
$ws.Cells.Item(3, 25).Value = 0.004740672436078721
$ws.Cells.Item(4, 3).Value = 0.3750152587890625
$ws.Cells.Item(4, 5).Value = 242.3591194254695
$ws.Cells.Item(4, 6).Value = 0.008034942080859841
$ws.Cells.Item(4, 7).Value = 0.006855174150669663
$ws.Cells.Item(4, 8).Value = 0.00637786138301151
$ws.Cells.Item(4, 9).Value = 0.006251224435706999
$ws.Cells.Item(4, 10).Value = 0.005779802418773928
$ws.Cells.Item(4, 11).Value = 0.005615643082687577
$ws.Cells.Item(4, 12).Value = 0.005553469060280929
$ws.Cells.Item(4, 13).Value = 0.005429948705476556
$ws.Cells.Item(4, 14).Value = 0.005371603514635828
$ws.Cells.Item(4, 15).Value = 0.00525042319174799
$ws.Cells.Item(4, 16).Value = 0.005061122399970292
$ws.Cells.Item(4, 17).Value = 0.004969103671680987
$ws.Cells.Item(4, 18).Value = 0.004950807933768671
$ws.Cells.Item(4, 19).Value = 0.004914848423611251
$ws.Cells.Item(4, 20).Value = 0.004898174143763114
$ws.Cells.Item(4, 21).Value = 0.004794656410900557
$ws.Cells.Item(4, 22).Value = 0.004751212617280174
$ws.Cells.Item(4, 23).Value = 0.004751212617280174
$ws.Cells.Item(4, 24).Value = 0.004732726138521165
$ws.Cells.Item(4, 25).Value = 0.004724349306539366
$ws.Cells.Item(5, 3).Value = 0.4851164817810059
$ws.Cells.Item(5, 5).Value = 239.3437704036314
$ws.Cells.Item(5, 6).Value = 0.008133045010140953
$ws.Cells.Item(5, 7).Value = 0.006639386608884319
$ws.Cells.Item(5, 8).Value = 0.005688398794631863
$ws.Cells.Item(5, 9).Value = 0.005688398794631863
$ws.Cells.Item(5, 10).Value = 0.005688398794631863
$ws.Cells.Item(5, 11).Value = 0.005527067381111986
$ws.Cells.Item(5, 12).Value = 0.005527067381111986
$ws.Cells.Item(5, 13).Value = 0.005327437206778037
$ws.Cells.Item(5, 14).Value = 0.005271873463982655
$ws.Cells.Item(5, 15).Value = 0.005193699501860978
$ws.Cells.Item(5, 16).Value = 0.00518870898396405
$ws.Cells.Item(5, 17).Value = 0.005083920554015383
$ws.Cells.Item(5, 18).Value = 0.00490652338490779
$ws.Cells.Item(5, 19).Value = 0.004861242829666563
$ws.Cells.Item(5, 20).Value = 0.004831561639743968
$ws.Cells.Item(5, 21).Value = 0.004794038677504409
$ws.Cells.Item(5, 22).Value = 0.00474511683968673
$ws.Cells.Item(5, 23).Value = 0.004707572323758403
$ws.Cells.Item(5, 24).Value = 0.004671818993633278
$ws.Cells.Item(5, 25).Value = 0.004665570573170202
$ws.Cells.Item(6, 3).Value = 0.3595192432403564
$ws.Cells.Item(6, 5).Value = 242.4663760195272
$ws.Cells.Item(6, 6).Value = 0.00816586172934892
$ws.Cells.Item(6, 7).Value = 0.006687867723263997
$ws.Cells.Item(6, 8).Value = 0.006362155963223285
$ws.Cells.Item(6, 9).Value = 0.006266715067506568
$ws.Cells.Item(6, 10).Value = 0.005902381103526024
$ws.Cells.Item(6, 11).Value = 0.005902381103526024
$ws.Cells.Item(6, 12).Value = 0.005197819504103579
$ws.Cells.Item(6, 13).Value = 0.005197819504103579
$ws.Cells.Item(6, 14).Value = 0.005197819504103579
$ws.Cells.Item(6, 15).Value = 0.005048946741824951
$ws.Cells.Item(6, 16).Value = 0.004987508421328176
$ws.Cells.Item(6, 17).Value = 0.004987508421328176
$ws.Cells.Item(6, 18).Value = 0.00498699478402533
$ws.Cells.Item(6, 19).Value = 0.004964430229732361
$ws.Cells.Item(6, 20).Value = 0.00492478077097446
$ws.Cells.Item(6, 21).Value = 0.004899229308722642
$ws.Cells.Item(6, 22).Value = 0.004809940817309257
$ws.Cells.Item(6, 23).Value = 0.004753425636822975
$ws.Cells.Item(6, 24).Value = 0.004726440078353354
$ws.Cells.Item(6, 25).Value = 0.004726440078353354
$ws.Cells.Item(7, 3).Value = 0.3906009197235107
$ws.Cells.Item(7, 5).Value = 240.9955944865906
$ws.Cells.Item(7, 6).Value = 0.008184182300238355
$ws.Cells.Item(7, 7).Value = 0.007234584370562798
$ws.Cells.Item(7, 8).Value = 0.006544502147698771
$ws.Cells.Item(7, 9).Value = 0.006286897058565648
$ws.Cells.Item(7, 10).Value = 0.005799016665886701
$ws.Cells.Item(7, 11).Value = 0.005799016665886701
$ws.Cells.Item(7, 12).Value = 0.005799016665886701
$ws.Cells.Item(7, 13).Value = 0.005615953751728501
$ws.Cells.Item(7, 14).Value = 0.005485799520823714
$ws.Cells.Item(7, 15).Value = 0.005369327052601477
$ws.Cells.Item(7, 16).Value = 0.0051517472257698
$ws.Cells.Item(7, 17).Value = 0.004892962476787694
$ws.Cells.Item(7, 18).Value = 0.004892962476787694
$ws.Cells.Item(7, 19).Value = 0.004892962476787694
$ws.Cells.Item(7, 20).Value = 0.004752660867629155
$ws.Cells.Item(7, 21).Value = 0.004752660867629155
$ws.Cells.Item(7, 22).Value = 0.004721567659416479
$ws.Cells.Item(7, 23).Value = 0.004721567659416479
$ws.Cells.Item(7, 24).Value = 0.004698011542235934
$ws.Cells.Item(7, 25).Value = 0.004697769873032953
$ws.Cells.Item(8, 3).Value = 0.375037670135498
$ws.Cells.Item(8, 5).Value = 247.0097872746501
$ws.Cells.Item(8, 6).Value = 0.008052843894278601
$ws.Cells.Item(8, 7).Value = 0.006870566186579535
$ws.Cells.Item(8, 8).Value = 0.005964125860644779
$ws.Cells.Item(8, 9).Value = 0.005616461236999806
$ws.Cells.Item(8, 10).Value = 0.005616461236999806
$ws.Cells.Item(8, 11).Value = 0.005516914915963665
$ws.Cells.Item(8, 12).Value = 0.005516914915963665
$ws.Cells.Item(8, 13).Value = 0.005426002392460753
$ws.Cells.Item(8, 14).Value = 0.005426002392460753
$ws.Cells.Item(8, 15).Value = 0.005293750545715702
$ws.Cells.Item(8, 16).Value = 0.005282849900755415
$ws.Cells.Item(8, 17).Value = 0.005141507275631957
$ws.Cells.Item(8, 18).Value = 0.005049141290239549
$ws.Cells.Item(8, 19).Value = 0.004996897224102059
$ws.Cells.Item(8, 20).Value = 0.004996897224102059
$ws.Cells.Item(8, 21).Value = 0.004960033490837797
$ws.Cells.Item(8, 22).Value = 0.004896057374006767
$ws.Cells.Item(8, 23).Value = 0.004873931665982855
$ws.Cells.Item(8, 24).Value = 0.00482079851533992
$ws.Cells.Item(8, 25).Value = 0.004815005599895712
$ws.Cells.Item(9, 3).Value = 0.4374630451202393
$ws.Cells.Item(9, 5).Value = 238.1654002635314
$ws.Cells.Item(9, 6).Value = 0.008002071430929755
$ws.Cells.Item(9, 7).Value = 0.006774340716937335
$ws.Cells.Item(9, 8).Value = 0.005937020190302925
$ws.Cells.Item(9, 9).Value = 0.005937020190302925
$ws.Cells.Item(9, 10).Value = 0.005937020190302925
$ws.Cells.Item(9, 11).Value = 0.005530727479320784
$ws.Cells.Item(9, 12).Value = 0.005291225092628611
$ws.Cells.Item(9, 13).Value = 0.005086456959976542
$ws.Cells.Item(9, 14).Value = 0.005086456959976542
$ws.Cells.Item(9, 15).Value = 0.005086456959976542
$ws.Cells.Item(9, 16).Value = 0.005075980086399603
$ws.Cells.Item(9, 17).Value = 0.004914367461619457
$ws.Cells.Item(9, 18).Value = 0.004908741667936125
$ws.Cells.Item(9, 19).Value = 0.004819093201887816
$ws.Cells.Item(9, 20).Value = 0.00478253231579279
$ws.Cells.Item(9, 21).Value = 0.004730070007545897
$ws.Cells.Item(9, 22).Value = 0.00472692041430743
$ws.Cells.Item(9, 23).Value = 0.004678645825660935
$ws.Cells.Item(9, 24).Value = 0.004654335157475679
$ws.Cells.Item(9, 25).Value = 0.004642600395000611
$ws.Cells.Item(10, 3).Value = 0.578150749206543
$ws.Cells.Item(10, 5).Value = 242.3299165268945
$ws.Cells.Item(10, 6).Value = 0.008537041040833346
$ws.Cells.Item(10, 7).Value = 0.006754188838317523
$ws.Cells.Item(10, 8).Value = 0.005978629942973778
$ws.Cells.Item(10, 9).Value = 0.005875066728381533
$ws.Cells.Item(10, 10).Value = 0.005678723423180558
$ws.Cells.Item(10, 11).Value = 0.005486174977076796
$ws.Cells.Item(10, 12).Value = 0.005435576610622861
$ws.Cells.Item(10, 13).Value = 0.005414632404476274
$ws.Cells.Item(10, 14).Value = 0.005344385784891604
$ws.Cells.Item(10, 15).Value = 0.005107680234120164
$ws.Cells.Item(10, 16).Value = 0.005061248014491872
$ws.Cells.Item(10, 17).Value = 0.005040686011735978
$ws.Cells.Item(10, 18).Value = 0.004974121204686099
$ws.Cells.Item(10, 19).Value = 0.004971773064760526
$ws.Cells.Item(10, 20).Value = 0.004853685569947312
$ws.Cells.Item(10, 21).Value = 0.004803575652806197
$ws.Cells.Item(10, 22).Value = 0.004786118641304089
$ws.Cells.Item(10, 23).Value = 0.004747052430725833
$ws.Cells.Item(10, 24).Value = 0.004747052430725833
$ws.Cells.Item(10, 25).Value = 0.004723780049257202
$ws.Cells.Item(11, 3).Value = 0.4062588214874268
$ws.Cells.Item(11, 5).Value = 239.9601500253921
$ws.Cells.Item(11, 6).Value = 0.008276527517561002
$ws.Cells.Item(11, 7).Value = 0.006981194037130389
$ws.Cells.Item(11, 8).Value = 0.006321667997540679
$ws.Cells.Item(11, 9).Value = 0.006098897532874562
$ws.Cells.Item(11, 10).Value = 0.005700851469481202
$ws.Cells.Item(11, 11).Value = 0.005507073656197064
$ws.Cells.Item(11, 12).Value = 0.005373432844486482
$ws.Cells.Item(11, 13).Value = 0.005373432844486482
$ws.Cells.Item(11, 14).Value = 0.00536664176229867
$ws.Cells.Item(11, 15).Value = 0.005359028083044446
$ws.Cells.Item(11, 16).Value = 0.005225607691886818
$ws.Cells.Item(11, 17).Value = 0.005209519000942576
$ws.Cells.Item(11, 18).Value = 0.005064045562738466
$ws.Cells.Item(11, 19).Value = 0.004852148839826481
$ws.Cells.Item(11, 20).Value = 0.004852148839826481
$ws.Cells.Item(11, 21).Value = 0.00477596950239993
$ws.Cells.Item(11, 22).Value = 0.004728757774879893
$ws.Cells.Item(11, 23).Value = 0.004728757774879893
$ws.Cells.Item(11, 24).Value = 0.004677585770475478
$ws.Cells.Item(11, 25).Value = 0.004677585770475478
